# FS-014: update the bad-address-test fixture with a new sample row.
#
# Row 2 held a single "Member ID" test value; it is replaced with a full
# sample record (Member ID / Postal Code / Unit Number / My Mailbox / PO Box)
# matching the Sheet1 header order: A=Member ID, B=Postal Code,
# C=Unit Number, D=My Mailbox, E=PO Box.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new sample values. C/D/E (new shared strings) and the numeric
# postal code are set before A2 so the shared-string table is appended in
# the same order the new strings appear in the saved workbook
# (200A, SGPost, 102A, then A300000067).
$ws.Range("C2").Value = "200A"
$ws.Range("D2").Value = "SGPost"
$ws.Range("E2").Value = "102A"
$ws.Range("A2").Value = "A300000067"
$ws.Range("B2").Value = 569933

# The postal-code cell (B2) picks up a left/top aligned style.
# -4131 = xlLeft, -4160 = xlTop
$ws.Range("B2").HorizontalAlignment = -4131
$ws.Range("B2").VerticalAlignment = -4160

# Column widths were hand-adjusted: A and B narrower, C/D/E back to the
# sheet's default width, and a new width set on (then-empty) column F.
$ws.Columns("A").ColumnWidth = 10
$ws.Columns("B").ColumnWidth = 8.67
$ws.Columns("F").ColumnWidth = 15.67

# The active selection moved to I8.
[void]$ws.Range("I8").Select()
